# Apply "bug fixes / new data" edit: insert four new teams (Delft, Leiden,
# Maastricht, Wageningen) into the participating-teams list, keeping the
# list sorted alphabetically by team name, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final alphabetically-sorted team list (header + 11 teams = rows 1..12)
$data = @(
    @("Teamnaam", "Plaats", "Inschrijflijst"),
    @("4 happy feet", "Enschede", "NTDS_Enschede.xlsx"),
    @("AmsterDance", "Amsterdam", "NTDS_Amsterdam.xlsx"),
    @("Blue Suede Shoes", "Delft", "NTDS_Delft.xlsx"),
    @("Dance Fever", "Nijmegen", "NTDS_Nijmegen.xlsx"),
    @("Erasmus Dance Society", "Rotterdam", "NTDS_Rotterdam.xlsx"),
    @("Footloose", "Eindhoven", "NTDS_Eindhoven.xlsx"),
    @("LeiDance", "Leiden", "NTDS_Leiden.xlsx"),
    @("Let's Dance", "Maastricht", "NTDS_Maastricht.xlsx"),
    @("The Blue Toes", "Groningen", "NTDS_Groningen.xlsx"),
    @("U Dance", "Utrecht", "NTDS_Utrecht.xlsx"),
    @("WUBDA", "Wageningen", "NTDS_Wageningen.xlsx")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Move/select the cell the editor left active in the source file.
$ws.Range("D5").Select()
